# Update "想去人数" (wishlist / "want to go" count) values in column F
# across the "展览", "演出" and "全部类型" sheets to reflect newly
# generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 1032    # was 1031
$ws.Range("F8").Value  = 195     # was 194
$ws.Range("F10").Value = 5       # was 4
$ws.Range("F14").Value = 12305   # was 12292
$ws.Range("F15").Value = 78      # was 74
$ws.Range("F16").Value = 5466    # was 5464

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value  = 114     # was 113

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 114     # was 113
$ws.Range("F7").Value  = 1032    # was 1031
$ws.Range("F10").Value = 195     # was 194
$ws.Range("F12").Value = 5       # was 4
$ws.Range("F16").Value = 12305   # was 12292
$ws.Range("F18").Value = 78      # was 74
$ws.Range("F19").Value = 5466    # was 5464
